# Auto-generated edit script applying the cryptos.xlsx data refresh
# (coin ranking reshuffle + updated Price/Volume(1h) figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.14%"

# Row 3
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.22%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.178"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.44%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05734"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.99%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.597"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.92%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8625"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.49%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8837"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.24%"

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.61%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07078"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.64%"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03272"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.88%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02874"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.73%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09403"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.19%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001516"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.46%"

# Row 15
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04148"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.25%"

# Row 16
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005995"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.82%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005980"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.11%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.505"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.87%"

# Row 19
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.069"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.63%"

# Row 20
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.180"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.08%"

# Row 21
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3184"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.30%"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.30%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.621"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.33%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001213"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.86%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004506"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.52%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001209"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "23.44%"

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.17%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03785"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.67%"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1072"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-22.02%"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002588"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.92%"

# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003519"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-41.91%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01005"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "23.20%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005114"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.90%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08893"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-18.37%"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.10%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
